$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin market data values, preserving text cell type
function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextCell "D2" "332.08"
Set-TextCell "E2" "1.81%"
Set-TextCell "G2" "3"
Set-TextCell "D3" "45.76"
Set-TextCell "E3" "3.29%"
Set-TextCell "G3" "3"
Set-TextCell "D4" "5.724"
Set-TextCell "E4" "4.14%"
Set-TextCell "G4" "3"
Set-TextCell "D5" "0.08361"
Set-TextCell "E5" "4.40%"
Set-TextCell "G5" "3"
Set-TextCell "D6" "2.051"
Set-TextCell "E6" "1.04%"
Set-TextCell "G6" "3"
Set-TextCell "D7" "0.9744"
Set-TextCell "E7" "2.76%"
Set-TextCell "G7" "3"
Set-TextCell "B8" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C8" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D8" "0.1157"
Set-TextCell "E8" "0.45%"
Set-TextCell "G8" "3"
Set-TextCell "B9" "WazirX"
Set-TextCell "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D9" "0.1942"
Set-TextCell "E9" "5.49%"
Set-TextCell "G9" "3"
Set-TextCell "B10" "MCDex"
Set-TextCell "C10" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D10" "10.44"
Set-TextCell "E10" "-14.08%"
Set-TextCell "G10" "3"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.09996"
Set-TextCell "E11" "2.40%"
Set-TextCell "G11" "3"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04630"
Set-TextCell "E12" "1.34%"
Set-TextCell "G12" "3"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.1059"
Set-TextCell "E13" "-0.53%"
Set-TextCell "G13" "3"
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001281"
Set-TextCell "E14" "1.36%"
Set-TextCell "G14" "3"
Set-TextCell "B15" "TigerCash"
Set-TextCell "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D15" "0.006061"
Set-TextCell "E15" "3.58%"
Set-TextCell "G15" "3"
Set-TextCell "B16" "LEO"
Set-TextCell "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D16" "3.371"
Set-TextCell "E16" "0.17%"
Set-TextCell "G16" "3"
Set-TextCell "B17" "GateToken"
Set-TextCell "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D17" "4.455"
Set-TextCell "E17" "3.75%"
Set-TextCell "G17" "3"
Set-TextCell "B18" "BTSEToken"
Set-TextCell "C18" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D18" "2.576"
Set-TextCell "E18" "0.18%"
Set-TextCell "G18" "3"
Set-TextCell "D19" "0.3350"
Set-TextCell "E19" "-3.71%"
Set-TextCell "G19" "3"
Set-TextCell "D20" "0.1401"
Set-TextCell "E20" "-0.31%"
Set-TextCell "G20" "3"
Set-TextCell "D21" "0.2595"
Set-TextCell "E21" "1.98%"
Set-TextCell "G21" "3"
Set-TextCell "D22" "0.04207"
Set-TextCell "E22" "3.15%"
Set-TextCell "G22" "3"
Set-TextCell "D23" "0.001308"
Set-TextCell "E23" "5.24%"
Set-TextCell "G23" "3"
Set-TextCell "D24" "0.004592"
Set-TextCell "E24" "6.71%"
Set-TextCell "G24" "3"
Set-TextCell "D25" "0.0001282"
Set-TextCell "G25" "3"
Set-TextCell "D26" "0.0003744"
Set-TextCell "E26" "0.07%"
Set-TextCell "G26" "3"
Set-TextCell "G27" "3"
Set-TextCell "G28" "3"
Set-TextCell "G29" "3"
Set-TextCell "G30" "3"
Set-TextCell "G31" "3"
Set-TextCell "G32" "3"
Set-TextCell "G33" "3"
Set-TextCell "G34" "3"
Set-TextCell "G35" "3"
Set-TextCell "G36" "3"
Set-TextCell "G37" "3"
Set-TextCell "D38" "0.02756"
Set-TextCell "E38" "7.45%"
Set-TextCell "G38" "3"
Set-TextCell "D39" "0.05844"
Set-TextCell "E39" "5.57%"
Set-TextCell "G39" "3"
Set-TextCell "D40" "0.007736"
Set-TextCell "E40" "2.81%"
Set-TextCell "G40" "3"
Set-TextCell "D41" "0.1436"
Set-TextCell "E41" "3.24%"
Set-TextCell "G41" "3"
Set-TextCell "D42" "0.007195"
Set-TextCell "E42" "-5.10%"
Set-TextCell "G42" "3"
Set-TextCell "D43" "0.001976"
Set-TextCell "E43" "-1.89%"
Set-TextCell "G43" "3"
Set-TextCell "D44" "0.008191"
Set-TextCell "E44" "-3.24%"
Set-TextCell "G44" "3"
Set-TextCell "D45" "0.00007200"
Set-TextCell "E45" "1.20%"
Set-TextCell "G45" "3"
Set-TextCell "D46" "0.00000000751"
Set-TextCell "E46" "0.18%"
Set-TextCell "G46" "3"
Set-TextCell "D47" "0.0005808"
Set-TextCell "E47" "-0.06%"
Set-TextCell "G47" "3"
Set-TextCell "D48" "0.003491"
Set-TextCell "E48" "-1.31%"
Set-TextCell "G48" "3"
Set-TextCell "D49" "0.003502"
Set-TextCell "E49" "52.15%"
Set-TextCell "G49" "3"
Set-TextCell "D50" "0.00002103"
Set-TextCell "E50" "0.18%"
Set-TextCell "G50" "3"
Set-TextCell "D51" "0.0002002"
Set-TextCell "E51" "0.18%"
Set-TextCell "G51" "3"
